# This document stores its photographs as legacy VML pictures
# (<w:pict><v:shape .../></w:pict>) rather than modern DrawingML, so they
# are not exposed through $d.Shapes / $d.InlineShapes in this runtime.
# The only reliable way to resize them is to rebuild the whole paragraph
# that hosts the <v:shape> with a corrected "width:...; height:...;" style
# and drop it back in with Range.InsertXML (which *replaces* the target
# range's contents).
#
# Each picture paragraph is immediately followed by a short, unique
# caption ("Frente", "Verso", ...). Find.Execute reliably reports the
# Start/End of that caption text, and the picture paragraph occupies the
# two characters (paragraph mark + cell mark) right before it, so we use
# that as the anchor instead of Paragraphs.Item(), whose cached
# Start/End can be stale for paragraphs that hold only a picture.

$d = $word.ActiveDocument

$pics = @(
    @{ Label = "Frente";                                 RId = "rId7";  Width = "220pt"; Height = "178.75pt" },
    @{ Label = "Verso";                                  RId = "rId8";  Width = "220pt"; Height = "188.65pt" },
    @{ Label = "Vista lateral direita";                  RId = "rId9";  Width = "250pt"; Height = "189.375pt" },
    @{ Label = "Vista lateral esquerda";                 RId = "rId10"; Width = "150pt"; Height = "111pt" },
    @{ Label = "Número de série";                        RId = "rId11"; Width = "150pt"; Height = "97.5pt" },
    @{ Label = "Base-Cartucho(s) calibre .380 AUTO";     RId = "rId12"; Width = "220pt"; Height = "164.45pt" },
    @{ Label = "Lateral-Cartucho(s) calibre .380 AUTO";  RId = "rId13"; Width = "220pt"; Height = "167.2pt" }
)

foreach ($pic in $pics) {
    $found = $d.Content.Find.Execute($pic.Label, $true, $false, $false, $false, $false,
                                      $true, 1, $false, "", 0)

    $label = $d.Content.Duplicate
    $label.Find.Execute($pic.Label, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

    $picRange = $d.Range($label.Start - 2, $label.Start)

    $style = "width:" + $pic.Width + "; height:" + $pic.Height + "; margin-left:0pt; margin-top:0pt; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line;"

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:v="urn:schemas-microsoft-com:vml" ' +
           'xmlns:w10="urn:schemas-microsoft-com:office:word" ' +
           'xmlns:o="urn:schemas-microsoft-com:office:office" ' +
           'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
           '<w:pPr><w:jc w:val="center"/></w:pPr>' +
           '<w:r><w:pict>' +
           '<v:shape type="#_x0000_t75" style="' + $style + '">' +
           '<w10:wrap type="inline"/>' +
           '<v:imagedata r:id="' + $pic.RId + '" o:title=""/>' +
           '</v:shape>' +
           '</w:pict></w:r>' +
           '</w:p>'

    $picRange.InsertXML($xml)
}
